$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.824.18"
$ws.Range("E2").Value = "  +4.18%  "
$ws.Range("D3").Value = "2.266.60"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.80"
$ws.Range("E5").Value = "  +3.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.70"
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("E7").Value = "  +3.59%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.43"
$ws.Range("E10").Value = "  +5.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.61"
$ws.Range("E11").Value = "  +3.76%  "
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").Value = "2.614.98"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.19"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "2.272.52"
$ws.Range("E17").Value = "  +4.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.762"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("D19").Value = "41.732.01"
$ws.Range("E19").Value = "  +4.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.24"
$ws.Range("E20").Value = "  +8.46%  "
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("E22").Value = "  +2.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.80"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "241.31"
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("E25").Value = "  +4.32%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +5.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.11"
$ws.Range("E28").Value = "  +3.61%  "
$ws.Range("E29").Value = "  +10.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.54"
$ws.Range("E30").Value = "  +2.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.28"
$ws.Range("E31").Value = "  +7.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "159.61"
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0745"
$ws.Range("E35").Value = "  +4.18%  "
$ws.Range("E36").Value = "  -2.44%  "
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("E38").Value = "  +7.33%  "
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("E40").Value = "  +3.89%  "
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.91"
$ws.Range("E42").Value = "  +3.73%  "
$ws.Range("D43").Value = "2.061.92"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.31"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.15"
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("E48").Value = "  +6.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "73.41"
$ws.Range("E49").Value = "  +8.08%  "
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.16"
$ws.Range("E51").Value = "  +2.49%  "
